$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("A3").Value = 70275566

# I3 holds a numeric-looking value but must stay text, like the rest of the column.
# A leading apostrophe forces Excel to store it as text instead of auto-converting
# it to a number, matching the original inlineStr text storage.
$ws.Range("I3").Value = "'146"

$ws.Range("K3").Value = "överblommad"
$ws.Range("P3").Value = "St Lommarstorp, Srm"
$ws.Range("S3").Value = 10
$ws.Range("X3").Value = "D-Str-0270"

# Y3/AA3 hold date-looking text values but must stay text, like the rest of the column.
$ws.Range("Y3").Value = "'2016-09-07"
$ws.Range("AA3").Value = "'2016-09-07"

$ws.Range("AH3").Value = "Skogsmark"
$ws.Range("AI3").Value = "Barrblandskog"
$ws.Range("AW3").Value = "Bo Karlsson"
$ws.Range("AX3").Value = "Bernt Andersson"
$ws.Range("AY3").Value = "Floraväkteri Sverige"
